{"js": "// Turn the \"Responses\" doc into a cover letter:\n//  - keep the existing first paragraph (it carries the _GoBack bookmark\n//    and the \"Reviewer Number\" style) but change its text from\n//    \"Reviewer: 1\" to \"Cover Letter\".\n//  - add a new normal-style paragraph right after it with the\n//    introductory sentence.\n//  - add a new \"Reviewer Number\" styled paragraph right after that one\n//    which now holds the \"Reviewer: 1\" heading that used to be first.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstPara = paragraphs.items[0];\nfirstPara.load(\"style\");\nawait context.sync();\n\nconst reviewerStyle = firstPara.style; // \"Reviewer Number\"\n\n// Insert the new \"Reviewer: 1\" heading paragraph right after the first\n// paragraph, using the same paragraph style as the original heading.\nconst reviewerPara = firstPara.insertParagraph(\"Reviewer: 1\", Word.InsertLocation.after);\nreviewerPara.style = reviewerStyle;\n\n// Insert the cover-letter intro paragraph between the first paragraph\n// and the new \"Reviewer: 1\" paragraph, using the default body style.\nconst introPara = firstPara.insertParagraph(\n  \"In this revision I have endeavored to address all of the reviewers\\u2019 comments.  Below is my response to each comment and a description of the changes made.\",\n  Word.InsertLocation.after\n);\nintroPara.style = \"Normal\";\n\n// Replace the original \"Reviewer: 1\" run text (inside the first\n// paragraph) with \"Cover Letter\", leaving the bookmark untouched.\nconst found = firstPara.search(\"Reviewer: 1\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nfound.items[0].insertText(\"Cover Letter\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$p1 = $d.Paragraphs(1)\n\n# Insert a new paragraph mark right after paragraph 1; this will hold\n# the \"Reviewer: 1\" heading that is moving out of paragraph 1.\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs(2)\n$p2.Range.Text = \"Reviewer: 1\"\n$p2.Style = $p1.Style\n\n# Insert another new paragraph mark after paragraph 1 (it lands between\n# paragraph 1 and the \"Reviewer: 1\" paragraph) for the cover-letter\n# intro sentence, using the document's default (Normal) style.\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs(2)\n$p2.Range.Text = \"In this revision I have endeavored to address all of the reviewers\" + [char]0x2019 + \" comments.  Below is my response to each comment and a description of the changes made.\"\n$p2.Style = \"Normal\"\n\n# Replace the original heading text in paragraph 1 with \"Cover Letter\",\n# keeping the paragraph's bookmark and formatting intact.\n$r = $p1.Range\n$textRange = $d.Range($r.Start, $r.End - 1)\n$textRange.Text = \"Cover Letter\"\n"}
